$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before current column F ("District") to hold the new "Address" column.
$ws.Columns.Item(6).Insert()

# Header
$ws.Range("F2").Value = "Address"

# Address values for each data row (3-50), derived from the school name
# embedded in column B (Names), matching the commit's new Address column.
$addresses = @{
    3  = "G H S Santpur"
    4  = "G P U C (HS) T K RoadBhadravathi"
    5  = "G J C (HS) Channagiri"
    6  = "Govt. P U College (H S) Aurad(B)"
    7  = "Sri Kaliveer High School Alur"
    8  = "S A G H S BirurKadur"
    9  = "Karnataka Public School BasavapatnaArakalgud"
    10 = "Shri Jayakeerthi High School Garag"
    11 = "Ambedkar High School Magge Alur"
    12 = "G H S DoddakanagaluAlur"
    13 = "G J C KallihalBhadravathi"
    14 = "S S T Girls H S Alnavar"
    15 = "Janatha High School BasawapattanaChannagiri"
    16 = "S G M High School Garag"
    17 = "S R H S GoppenahalliChannagiri"
    18 = "Shri. B S Patil P U College (HS) ManagaliBasvan Bagewadi"
    19 = "Govt. G H S Ilkal"
    20 = "Govt. H S Thanakushnoor"
    22 = "Katageri Vidyavardak Sangad High School KatageriBadami"
    23 = "Shri Basveshwar Govt. P U CollegeBasvan Bagewadi"
    24 = "Govt. Boys High School Wadagaon(D)Aurad"
    25 = "Shree Guru Virupaksheswara High School UppinBetagiri"
    26 = "G P U CollegeHolehonnurBhadravathi"
    27 = "Govt. Boys Composite PU College Channapatna"
    28 = "Shree Shanteshwar High School Amminabhavi"
    29 = "G G P U C New Town Bhadravathi"
    30 = "G J C Channagiri"
    31 = "S J M H S BirurKadur"
    32 = "S V H School S bidare"
    33 = "S C H S TanigebyluBirurTarikere"
    34 = "Shree Lakshmeesha High SchoolDevanurKadur"
    35 = "Veer Pulikeshi H S Badami"
    36 = "G H S Chawar DapkaAurad"
    37 = "Sri Annadaneshular Comp. Jr. CollegeBelurBadami"
    38 = "Smt. A C Ghattad Girls High School KerurBadami"
    39 = "G H S MashalAfzalpur"
    40 = "G H S KanathurAlur"
    41 = "S R High School IngaleswarBasavan Bagewadi"
    42 = "G H S SaganoorAfzalapur"
    43 = "G J C NallurChannagiri"
    44 = "Sanchi Honnamma G G P U C Old Town Bhadravathi"
    45 = "G H S KarajagiAfazalapur"
    46 = "Govt. Girls P U CollegeJamkhandi High Section"
    47 = "Govt. High School BalliganurBirur"
    48 = "R T Desai Govt. P U College NeerbudihalBadami"
    49 = "N E S Alnavar"
    50 = "G G H S Bilagi"
}

foreach ($row in ($addresses.Keys | Sort-Object)) {
    $ws.Range("F$row").Value = $addresses[$row]
}

# Row 21 has no address value (stays blank) but still needs the inlineStr cell type,
# which Excel will naturally keep as empty text since the column insert already
# shifted it from the old F21.
